$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '23.445.74'
$ws.Range("E2").Value = '  -1.08%  '
# Row 3
$ws.Range("D3").Value = '1.638.56'
$ws.Range("E3").Value = '  -0.96%  '
# Row 4
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.30%  '
# Row 5
$ws.Range("D5").Value = '''0.9994'
$ws.Range("E5").Value = '  +0.09%  '
# Row 6
$ws.Range("D6").Value = '''303.31'
$ws.Range("E6").Value = '  -1.04%  '
# Row 7
$ws.Range("D7").Value = '''0.3788'
$ws.Range("E7").Value = '  +0.35%  '
# Row 8
$ws.Range("E8").Value = '  -0.89%  '
# Row 9
$ws.Range("D9").Value = '''0.3630'
$ws.Range("E9").Value = '  -1.02%  '
# Row 10
$ws.Range("D10").Value = '''0.08194'
$ws.Range("E10").Value = '  +0.24%  '
# Row 11
$ws.Range("D11").Value = '''1.230'
$ws.Range("E11").Value = '  -3.56%  '
# Row 12
$ws.Range("D12").Value = '''1.001'
$ws.Range("E12").Value = '  +0.24%  '
# Row 13
$ws.Range("D13").Value = '''22.47'
$ws.Range("E13").Value = '  -2.98%  '
# Row 14
$ws.Range("D14").Value = '''6.484'
$ws.Range("E14").Value = '  -3.39%  '
# Row 15
$ws.Range("D15").Value = '''7.379'
$ws.Range("E15").Value = '  -0.13%  '
# Row 16
$ws.Range("D16").Value = '''0.00001240'
$ws.Range("E16").Value = '  -3.02%  '
# Row 17
$ws.Range("D17").Value = '1.631.91'
$ws.Range("E17").Value = '  -1.86%  '
# Row 18
$ws.Range("D18").Value = '''95.26'
$ws.Range("E18").Value = '  -0.25%  '
# Row 19
$ws.Range("D19").Value = '''0.06945'
$ws.Range("E19").Value = '  +0.43%  '
# Row 20
$ws.Range("D20").Value = '''6.598'
$ws.Range("E20").Value = '  -0.16%  '
# Row 21
$ws.Range("D21").Value = '''17.52'
$ws.Range("E21").Value = '  -5.05%  '
# Row 22
$ws.Range("D22").Value = '''1.002'
$ws.Range("E22").Value = '  +0.39%  '
# Row 23
$ws.Range("D23").Value = '''12.55'
$ws.Range("E23").Value = '  -3.72%  '
# Row 24
$ws.Range("D24").Value = '23.452.06'
$ws.Range("E24").Value = '  -1.03%  '
# Row 25
$ws.Range("D25").Value = '''2.519'
$ws.Range("E25").Value = '  +4.34%  '
# Row 26
$ws.Range("D26").Value = '''3.078'
$ws.Range("E26").Value = '  -3.14%  '
# Row 27
$ws.Range("D27").Value = '''21.20'
$ws.Range("E27").Value = '  -1.30%  '
# Row 28
$ws.Range("D28").Value = '''151.62'
$ws.Range("E28").Value = '  +0.25%  '
# Row 29
$ws.Range("D29").Value = '''5.265'
$ws.Range("E29").Value = '  -0.95%  '
# Row 30
$ws.Range("D30").Value = '''133.87'
$ws.Range("E30").Value = '  -2.40%  '
# Row 31
$ws.Range("D31").Value = '1.814.91'
$ws.Range("E31").Value = '  -1.81%  '
# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''6.639'
$ws.Range("E32").Value = '  -4.24%  '
# Row 33
$ws.Range("B33").Value = 'WEMIXTOKEN'
$ws.Range("C33").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D33").Value = '''2.155'
$ws.Range("E33").Value = '  -8.27%  '
# Row 34
$ws.Range("D34").Value = '''1.053'
$ws.Range("E34").Value = '  +7.88%  '
# Row 35
$ws.Range("D35").Value = '''11.35'
$ws.Range("E35").Value = '  +2.43%  '
# Row 36
$ws.Range("D36").Value = '''0.02759'
$ws.Range("E36").Value = '  -4.43%  '
# Row 37
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").Value = '''0.2497'
$ws.Range("E37").Value = '  -3.47%  '
# Row 38
$ws.Range("B38").Value = 'Stellar'
$ws.Range("C38").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D38").Value = '''0.08784'
$ws.Range("E38").Value = '  -1.53%  '
# Row 39
$ws.Range("D39").Value = '''0.07126'
$ws.Range("E39").Value = '  -3.39%  '
# Row 40
$ws.Range("D40").Value = '''6.047'
$ws.Range("E40").Value = '  -5.54%  '
# Row 41
$ws.Range("D41").Value = '''0.7032'
$ws.Range("E41").Value = '  -2.79%  '
# Row 42
$ws.Range("D42").Value = '''1.340'
$ws.Range("E42").Value = '  -3.11%  '
# Row 43
$ws.Range("D43").Value = '''12.23'
$ws.Range("E43").Value = '  -4.00%  '
# Row 44
$ws.Range("D44").Value = '''15.84'
$ws.Range("E44").Value = '  -4.27%  '
# Row 45
$ws.Range("D45").Value = '''0.6523'
$ws.Range("E45").Value = '  -2.03%  '
# Row 46
$ws.Range("D46").Value = '''0.9991'
$ws.Range("E46").Value = '  +0.22%  '
# Row 47
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '''2.279'
$ws.Range("E47").Value = '  -4.27%  '
# Row 48
$ws.Range("B48").Value = 'PancakeSwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D48").Value = '''3.969'
$ws.Range("E48").Value = '  -1.47%  '
# Row 49
$ws.Range("D49").Value = '''0.08006'
$ws.Range("E49").Value = '  -0.41%  '
# Row 50
$ws.Range("D50").Value = '''127.36'
$ws.Range("E50").Value = '  -0.83%  '
# Row 51
$ws.Range("D51").Value = '''1.195'
$ws.Range("E51").Value = '  -2.79%  '
